# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (Total) sheet
#    and populate it with the per-fund holding detail for 2022-Q1.
# 2. Insert a new summary row at the top of the "总计" sheet's data
#    (right after the header row) for the "2022-Q1" quarter, shifting the
#    previously-existing rows down and renumbering the index column (A).

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet, inserted immediately before "总计"
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Header row (row 1) - same style as the other detail sheets (bold / boxed)
$q1Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$q1HeaderCols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $q1Headers.Length; $i++) {
    $cell = $q1.Range($q1HeaderCols[$i] + "1")
    $cell.Value = $q1Headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Detail rows - B..G are textual (as in the other quarter sheets), H is numeric
$q1Rows = @(
    @("870009", "广发资管平衡精选一年持有混合A", "11.34", "94.29", "4.07", "0.4615", 10),
    @("872019", "广发资管平衡精选一年持有混合C", "1.54", "94.29", "4.07", "0.0627", 10)
)

for ($r = 0; $r -lt $q1Rows.Length; $r++) {
    $rowNum = $r + 2
    $q1.Range("A" + $rowNum).Value = $r

    $q1.Range("B" + $rowNum).Value = "'" + $q1Rows[$r][0]
    $q1.Range("C" + $rowNum).Value = $q1Rows[$r][1]
    $q1.Range("D" + $rowNum).Value = "'" + $q1Rows[$r][2]
    $q1.Range("E" + $rowNum).Value = "'" + $q1Rows[$r][3]
    $q1.Range("F" + $rowNum).Value = "'" + $q1Rows[$r][4]
    $q1.Range("G" + $rowNum).Value = "'" + $q1Rows[$r][5]
    $q1.Range("H" + $rowNum).Value = $q1Rows[$r][6]
}

# ---------------------------------------------------------------------
# 2. Prepend the 2022-Q1 summary row to the "总计" sheet
# ---------------------------------------------------------------------
# NOTE: re-fetch the sheet by name - the earlier $totalSheet handle now
# tracks whatever sheet sits at its old position (the freshly inserted
# "2022-Q1" sheet), since inserting a sheet shifted "总计" along.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The freshly-inserted blank row 2 doesn't carry the same per-column
# formatting as the data rows below it (row 3 is the old row 2, which
# still has the right look). Pull formats across before writing values,
# mirroring how "A" (the index column) is boxed/bold while B:D are plain.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.52

# Renumber the index column for the rows that shifted down
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
